# Fix oversized font in the "X" absence-marker cells of the Mittwoch (Wednesday)
# column of the attendance table: the cells were left at the default 11pt
# (no explicit run/paragraph-mark size) while all the other day-columns
# already carry an explicit 10pt (sz/szCs = 20 half-points) size.
#
# This sets the font size to 10pt on every cell in column 10 whose content
# is exactly "X" and which isn't already at 10pt, which reproduces the
# <w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr> addition to both the
# paragraph mark and the run, on each of the affected rows.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$rows = $t.Rows.Count
$targetCol = 10

for ($r = 1; $r -le $rows; $r++) {
    $cell = $t.Cell($r, $targetCol)
    $txt = $cell.Range.Text -replace "[\x00-\x1f]", ""
    if ($txt -eq "X") {
        if ($cell.Range.Font.Size -ne 10) {
            $cell.Range.Font.Size = 10
            $cell.Range.Font.SizeBi = 10
        }
    }
}
